$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Book Inventory" ---
$ws1 = $wb.Worksheets.Item("Book Inventory")

# Update The Hunger Games row: ISBN + quantities
$ws1.Range("C2").NumberFormat = "@"
$ws1.Range("C2").Value = "9780439023528"
$ws1.Range("D2").Value = 1
$ws1.Range("E2").Value = 1

# Replace row 3 (previously "Wild Things Storytelling Kit") with "The Outsiders"
$ws1.Range("A3").Value = "The Outsiders"
$ws1.Range("B3").Value = "S.E. Hinton"
$ws1.Range("C3").NumberFormat = "@"
$ws1.Range("C3").Value = "9780140385724"
$ws1.Range("D3").Value = 1
$ws1.Range("E3").Value = 1

# Add new row 4 for "To Kill A Mockingbird"
$ws1.Range("A4").Value = "To Kill A Mockingbird"
$ws1.Range("B4").Value = "by Harper Lee"
$ws1.Range("C4").NumberFormat = "@"
$ws1.Range("C4").Value = "9780446310789"
$ws1.Range("D4").Value = 2
$ws1.Range("E4").Value = 2

# --- Sheet 2: "Check Out-In" ---
# Remove the empty placeholder row that was used for testing
$ws2 = $wb.Worksheets.Item("Check Out-In")
$ws2.Rows.Item(3).Delete()
